$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint_One")

# Update "Story Points (Done that day)" entries that changed (E2 and G2 went from 0 to 1).
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 1

# Move the active cell selection to M9, matching the saved sheet view state.
$ws.Range("M9").Select()

$wb.Save()
